$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows before the current totals row (row 16) ---
# This pushes the old row 16 (totals) -> row 18 and old row 17 (footer) -> row 19,
# shifting the existing merged cells down automatically.
$ws.Rows("16:17").Insert()

# --- Clone the formatting of the last item row (row 15) onto the two new rows ---
$ws.Range("A15:Q15").Copy()
$ws.Range("A16:Q16").PasteSpecial(-4122)
$ws.Range("A15:Q15").Copy()
$ws.Range("A17:Q17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row heights for the two new item rows (totals/footer rows keep their own) ---
$ws.Rows("16").RowHeight = 25.5
$ws.Rows("17").RowHeight = 25.5
$ws.Rows("18").RowHeight = 24.75

# --- Re-create the merges inside the two new rows (A:B, C:G, H:K, L:M, N:O) ---
$ws.Range("A16:B16").Merge()
$ws.Range("C16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()
$ws.Range("N16:O16").Merge()

$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

# --- New row 16: item #10 ---
$ws.Range("A16").Value = 10
$ws.Range("C16").Value = "فرشاة اطفال ريتش ديلي"
$ws.Range("H16").Value = "26:0"
$ws.Range("L16").Value = "0"
$ws.Range("N16").Value = "15.00"
$ws.Range("P16").Value = "15.0000"
$ws.Range("Q16").Value = "1:0"

# --- New row 17: item #11 ---
$ws.Range("A17").Value = 11
$ws.Range("C17").Value = "فرشه اسنان POWER GOLD كبار"
$ws.Range("H17").Value = "11:0"
$ws.Range("L17").Value = "0"
$ws.Range("N17").Value = "15.00"
$ws.Range("P17").Value = "15.0000"
$ws.Range("Q17").Value = "1:0"

# --- Update the grand-total row (now row 18) ---
$ws.Range("P18").Value = 619.245

# --- Update the generated-at timestamp in the footer (now row 19) ---
$ws.Range("A19").Value = "Saturday, 27 September, 2025 10:55 AM"
